$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at D:E, shifting existing D:K to F:M
$ws.Range("D:E").Insert()

# Copy number formats into the two new columns so they match the rest of the table
$ws.Range("F8:G102").Copy()
$ws.Range("D8:E102").PasteSpecial(-4122)
$ws.Range("F7:G7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F38:G38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F80:G80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the final values for every data cell D:M, row by row
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(7, 6).Value = 43281
$ws.Cells.Item(7, 7).Value = 43190
$ws.Cells.Item(7, 8).Value = 43100
$ws.Cells.Item(7, 9).Value = 43008
$ws.Cells.Item(7, 10).Value = 42916
$ws.Cells.Item(7, 11).Value = 42825
$ws.Cells.Item(7, 12).Value = 42735
$ws.Cells.Item(7, 13).Value = 42643
$ws.Cells.Item(8, 4).Value = 2189000
$ws.Cells.Item(8, 5).Value = 2259000
$ws.Cells.Item(8, 6).Value = 2114000
$ws.Cells.Item(8, 7).Value = 1981000
$ws.Cells.Item(8, 8).Value = 1281400
$ws.Cells.Item(8, 9).Value = 2185000
$ws.Cells.Item(8, 10).Value = 2192000
$ws.Cells.Item(8, 11).Value = 2059000
$ws.Cells.Item(8, 12).Value = 8822000
$ws.Cells.Item(8, 13).Value = 2541200
$ws.Cells.Item(9, 4).Value = 108000
$ws.Cells.Item(9, 5).Value = 86000
$ws.Cells.Item(9, 6).Value = 74000
$ws.Cells.Item(9, 7).Value = 100000
$ws.Cells.Item(9, 8).Value = -34400
$ws.Cells.Item(9, 9).Value = 209000
$ws.Cells.Item(9, 10).Value = 195000
$ws.Cells.Item(9, 11).Value = 64000
$ws.Cells.Item(9, 12).Value = 415000
$ws.Cells.Item(9, 13).Value = 180700
$ws.Cells.Item(10, 4).Value = 2081000
$ws.Cells.Item(10, 5).Value = 2173000
$ws.Cells.Item(10, 6).Value = 2040000
$ws.Cells.Item(10, 7).Value = 1881000
$ws.Cells.Item(10, 8).Value = 1315800
$ws.Cells.Item(10, 9).Value = 1976000
$ws.Cells.Item(10, 10).Value = 1997000
$ws.Cells.Item(10, 11).Value = 1995000
$ws.Cells.Item(10, 12).Value = 8407000
$ws.Cells.Item(10, 13).Value = 2360500
$ws.Cells.Item(11, 4).ClearContents()
$ws.Cells.Item(11, 5).ClearContents()
$ws.Cells.Item(11, 6).ClearContents()
$ws.Cells.Item(11, 7).ClearContents()
$ws.Cells.Item(11, 8).ClearContents()
$ws.Cells.Item(11, 9).ClearContents()
$ws.Cells.Item(11, 10).ClearContents()
$ws.Cells.Item(11, 11).ClearContents()
$ws.Cells.Item(11, 12).ClearContents()
$ws.Cells.Item(11, 13).ClearContents()
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(12, 6).Value = "NA"
$ws.Cells.Item(12, 7).Value = "NA"
$ws.Cells.Item(12, 8).Value = "NA"
$ws.Cells.Item(12, 9).Value = "NA"
$ws.Cells.Item(12, 10).Value = "NA"
$ws.Cells.Item(12, 11).Value = "NA"
$ws.Cells.Item(12, 12).Value = "NA"
$ws.Cells.Item(12, 13).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(14, 4).Value = 7000
$ws.Cells.Item(14, 5).Value = 7000
$ws.Cells.Item(14, 6).Value = 1000
$ws.Cells.Item(14, 7).Value = 14000
$ws.Cells.Item(14, 8).Value = 15800
$ws.Cells.Item(14, 9).Value = 1000
$ws.Cells.Item(14, 10).Value = 1000
$ws.Cells.Item(14, 11).Value = 1000
$ws.Cells.Item(14, 12).Value = 34000
$ws.Cells.Item(14, 13).Value = 12700
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(16, 4).ClearContents()
$ws.Cells.Item(16, 5).ClearContents()
$ws.Cells.Item(16, 6).ClearContents()
$ws.Cells.Item(16, 7).ClearContents()
$ws.Cells.Item(16, 8).ClearContents()
$ws.Cells.Item(16, 9).ClearContents()
$ws.Cells.Item(16, 10).ClearContents()
$ws.Cells.Item(16, 11).ClearContents()
$ws.Cells.Item(16, 12).ClearContents()
$ws.Cells.Item(16, 13).ClearContents()
$ws.Cells.Item(17, 4).Value = 1937000
$ws.Cells.Item(17, 5).Value = 2018000
$ws.Cells.Item(17, 6).Value = 1811000
$ws.Cells.Item(17, 7).Value = 1905000
$ws.Cells.Item(17, 8).Value = 1000900
$ws.Cells.Item(17, 9).Value = 2078000
$ws.Cells.Item(17, 10).Value = 1965000
$ws.Cells.Item(17, 11).Value = 1883000
$ws.Cells.Item(17, 12).Value = 8422000
$ws.Cells.Item(17, 13).Value = 2825000
$ws.Cells.Item(18, 4).Value = 252000
$ws.Cells.Item(18, 5).Value = 241000
$ws.Cells.Item(18, 6).Value = 303000
$ws.Cells.Item(18, 7).Value = 76000
$ws.Cells.Item(18, 8).Value = 280500
$ws.Cells.Item(18, 9).Value = 107000
$ws.Cells.Item(18, 10).Value = 227000
$ws.Cells.Item(18, 11).Value = 176000
$ws.Cells.Item(18, 12).Value = 400000
$ws.Cells.Item(18, 13).Value = -283800
$ws.Cells.Item(19, 4).ClearContents()
$ws.Cells.Item(19, 5).ClearContents()
$ws.Cells.Item(19, 6).ClearContents()
$ws.Cells.Item(19, 7).ClearContents()
$ws.Cells.Item(19, 8).ClearContents()
$ws.Cells.Item(19, 9).ClearContents()
$ws.Cells.Item(19, 10).ClearContents()
$ws.Cells.Item(19, 11).ClearContents()
$ws.Cells.Item(19, 12).ClearContents()
$ws.Cells.Item(19, 13).ClearContents()
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = 0
$ws.Cells.Item(21, 4).Value = "NA"
$ws.Cells.Item(21, 5).Value = "NA"
$ws.Cells.Item(21, 6).Value = "NA"
$ws.Cells.Item(21, 7).Value = "NA"
$ws.Cells.Item(21, 8).Value = "NA"
$ws.Cells.Item(21, 9).Value = "NA"
$ws.Cells.Item(21, 10).Value = "NA"
$ws.Cells.Item(21, 11).Value = "NA"
$ws.Cells.Item(21, 12).Value = "NA"
$ws.Cells.Item(21, 13).Value = "NA"
$ws.Cells.Item(22, 4).Value = 90000
$ws.Cells.Item(22, 5).Value = 55000
$ws.Cells.Item(22, 6).Value = 62000
$ws.Cells.Item(22, 7).Value = 55000
$ws.Cells.Item(22, 8).Value = 62300
$ws.Cells.Item(22, 9).Value = 67000
$ws.Cells.Item(22, 10).Value = 72000
$ws.Cells.Item(22, 11).Value = 63000
$ws.Cells.Item(22, 12).Value = 390000
$ws.Cells.Item(22, 13).Value = 72100
$ws.Cells.Item(23, 4).Value = 162000
$ws.Cells.Item(23, 5).Value = 186000
$ws.Cells.Item(23, 6).Value = 241000
$ws.Cells.Item(23, 7).Value = 21000
$ws.Cells.Item(23, 8).Value = 218200
$ws.Cells.Item(23, 9).Value = 40000
$ws.Cells.Item(23, 10).Value = 155000
$ws.Cells.Item(23, 11).Value = 113000
$ws.Cells.Item(23, 12).Value = 10000
$ws.Cells.Item(23, 13).Value = -355900
$ws.Cells.Item(24, 4).Value = -23000
$ws.Cells.Item(24, 5).Value = 21000
$ws.Cells.Item(24, 6).Value = 45000
$ws.Cells.Item(24, 7).Value = 4000
$ws.Cells.Item(24, 8).Value = 42000
$ws.Cells.Item(24, 9).Value = -40000
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 93000
$ws.Cells.Item(24, 12).Value = -29000
$ws.Cells.Item(24, 13).Value = -119400
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(26, 4).Value = 185000
$ws.Cells.Item(26, 5).Value = 165000
$ws.Cells.Item(26, 6).Value = 196000
$ws.Cells.Item(26, 7).Value = 17000
$ws.Cells.Item(26, 8).Value = 176200
$ws.Cells.Item(26, 9).Value = 80000
$ws.Cells.Item(26, 10).Value = 155000
$ws.Cells.Item(26, 11).Value = 20000
$ws.Cells.Item(26, 12).Value = 39000
$ws.Cells.Item(26, 13).Value = -236500
$ws.Cells.Item(27, 4).Value = 129000
$ws.Cells.Item(27, 5).Value = 142000
$ws.Cells.Item(27, 6).Value = 138000
$ws.Cells.Item(27, 7).Value = 17000
$ws.Cells.Item(27, 8).Value = 94700
$ws.Cells.Item(27, 9).Value = 15000
$ws.Cells.Item(27, 10).Value = 103000
$ws.Cells.Item(27, 11).Value = 19000
$ws.Cells.Item(27, 12).Value = 10000
$ws.Cells.Item(27, 13).Value = -248100
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0
$ws.Cells.Item(29, 4).Value = -8000
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 28000
$ws.Cells.Item(29, 7).Value = 429000
$ws.Cells.Item(29, 8).Value = -3259000
$ws.Cells.Item(29, 9).Value = 134000
$ws.Cells.Item(29, 10).Value = 64000
$ws.Cells.Item(29, 11).Value = -162000
$ws.Cells.Item(29, 12).Value = -337000
$ws.Cells.Item(29, 13).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 0
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = 0
$ws.Cells.Item(33, 4).Value = 121000
$ws.Cells.Item(33, 5).Value = 142000
$ws.Cells.Item(33, 6).Value = 166000
$ws.Cells.Item(33, 7).Value = 446000
$ws.Cells.Item(33, 8).Value = -3164300
$ws.Cells.Item(33, 9).Value = 149000
$ws.Cells.Item(33, 10).Value = 167000
$ws.Cells.Item(33, 11).Value = -143000
$ws.Cells.Item(33, 12).Value = -327000
$ws.Cells.Item(33, 13).Value = -248100
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 0
$ws.Cells.Item(35, 4).Value = 121000
$ws.Cells.Item(35, 5).Value = 142000
$ws.Cells.Item(35, 6).Value = 166000
$ws.Cells.Item(35, 7).Value = 446000
$ws.Cells.Item(35, 8).Value = -3164300
$ws.Cells.Item(35, 9).Value = 149000
$ws.Cells.Item(35, 10).Value = 167000
$ws.Cells.Item(35, 11).Value = -143000
$ws.Cells.Item(35, 12).Value = -327000
$ws.Cells.Item(35, 13).Value = -248100
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(38, 6).Value = 43281
$ws.Cells.Item(38, 7).Value = 43190
$ws.Cells.Item(38, 8).Value = 43100
$ws.Cells.Item(38, 9).Value = 43008
$ws.Cells.Item(38, 10).Value = 42916
$ws.Cells.Item(38, 11).Value = 42825
$ws.Cells.Item(38, 12).Value = 42735
$ws.Cells.Item(38, 13).Value = 42643
$ws.Cells.Item(39, 4).ClearContents()
$ws.Cells.Item(39, 5).ClearContents()
$ws.Cells.Item(39, 6).ClearContents()
$ws.Cells.Item(39, 7).ClearContents()
$ws.Cells.Item(39, 8).ClearContents()
$ws.Cells.Item(39, 9).ClearContents()
$ws.Cells.Item(39, 10).ClearContents()
$ws.Cells.Item(39, 11).ClearContents()
$ws.Cells.Item(39, 12).ClearContents()
$ws.Cells.Item(39, 13).ClearContents()
$ws.Cells.Item(40, 4).ClearContents()
$ws.Cells.Item(40, 5).ClearContents()
$ws.Cells.Item(40, 6).ClearContents()
$ws.Cells.Item(40, 7).ClearContents()
$ws.Cells.Item(40, 8).ClearContents()
$ws.Cells.Item(40, 9).ClearContents()
$ws.Cells.Item(40, 10).ClearContents()
$ws.Cells.Item(40, 11).ClearContents()
$ws.Cells.Item(40, 12).ClearContents()
$ws.Cells.Item(40, 13).ClearContents()
$ws.Cells.Item(41, 4).Value = 1869000
$ws.Cells.Item(41, 5).Value = 1885000
$ws.Cells.Item(41, 6).Value = 1620000
$ws.Cells.Item(41, 7).Value = 1597000
$ws.Cells.Item(41, 8).Value = 1435000
$ws.Cells.Item(41, 9).Value = 2072900
$ws.Cells.Item(41, 10).Value = 2713300
$ws.Cells.Item(41, 11).Value = 2531300
$ws.Cells.Item(41, 12).Value = 2229000
$ws.Cells.Item(41, 13).Value = 2951600
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = 0
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(43, 5).Value = 0
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = 0
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).Value = 0
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = 0
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 0
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = 0
$ws.Cells.Item(47, 4).Value = 62012000
$ws.Cells.Item(47, 5).Value = 62277000
$ws.Cells.Item(47, 6).Value = 62241000
$ws.Cells.Item(47, 7).Value = 63014000
$ws.Cells.Item(47, 8).Value = 64644000
$ws.Cells.Item(47, 9).Value = 93681000
$ws.Cells.Item(47, 10).Value = 90575600
$ws.Cells.Item(47, 11).Value = 89116700
$ws.Cells.Item(47, 12).Value = 61153000
$ws.Cells.Item(47, 13).Value = 93051100
$ws.Cells.Item(48, 4).Value = "NA"
$ws.Cells.Item(48, 5).Value = "NA"
$ws.Cells.Item(48, 6).Value = "NA"
$ws.Cells.Item(48, 7).Value = "NA"
$ws.Cells.Item(48, 8).Value = 107000
$ws.Cells.Item(48, 9).Value = "NA"
$ws.Cells.Item(48, 10).Value = "NA"
$ws.Cells.Item(48, 11).Value = "NA"
$ws.Cells.Item(48, 12).Value = 112300
$ws.Cells.Item(48, 13).Value = "NA"
$ws.Cells.Item(49, 4).Value = "NA"
$ws.Cells.Item(49, 5).Value = "NA"
$ws.Cells.Item(49, 6).Value = "NA"
$ws.Cells.Item(49, 7).Value = "NA"
$ws.Cells.Item(49, 8).Value = 186000
$ws.Cells.Item(49, 9).Value = 196000
$ws.Cells.Item(49, 10).Value = 197800
$ws.Cells.Item(49, 11).Value = 208000
$ws.Cells.Item(49, 12).Value = 219500
$ws.Cells.Item(49, 13).Value = 227300
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = 0
$ws.Cells.Item(52, 4).Value = 1157000
$ws.Cells.Item(52, 5).Value = 1122000
$ws.Cells.Item(52, 6).Value = 1266000
$ws.Cells.Item(52, 7).Value = 58102000
$ws.Cells.Item(52, 8).Value = 59833000
$ws.Cells.Item(52, 9).Value = 1663700
$ws.Cells.Item(52, 10).Value = 1756900
$ws.Cells.Item(52, 11).Value = 2055700
$ws.Cells.Item(52, 12).Value = 66368800
$ws.Cells.Item(52, 13).Value = 1128100
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 0
$ws.Cells.Item(54, 4).Value = 154682000
$ws.Cells.Item(54, 5).Value = 165036000
$ws.Cells.Item(54, 6).Value = 162740000
$ws.Cells.Item(54, 7).Value = 219824000
$ws.Cells.Item(54, 8).Value = 222532000
$ws.Cells.Item(54, 9).Value = 226643900
$ws.Cells.Item(54, 10).Value = 219838800
$ws.Cells.Item(54, 11).Value = 217025500
$ws.Cells.Item(54, 12).Value = 214585000
$ws.Cells.Item(54, 13).Value = 218622200
$ws.Cells.Item(55, 4).ClearContents()
$ws.Cells.Item(55, 5).ClearContents()
$ws.Cells.Item(55, 6).ClearContents()
$ws.Cells.Item(55, 7).ClearContents()
$ws.Cells.Item(55, 8).ClearContents()
$ws.Cells.Item(55, 9).ClearContents()
$ws.Cells.Item(55, 10).ClearContents()
$ws.Cells.Item(55, 11).ClearContents()
$ws.Cells.Item(55, 12).ClearContents()
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(56, 4).ClearContents()
$ws.Cells.Item(56, 5).ClearContents()
$ws.Cells.Item(56, 6).ClearContents()
$ws.Cells.Item(56, 7).ClearContents()
$ws.Cells.Item(56, 8).ClearContents()
$ws.Cells.Item(56, 9).ClearContents()
$ws.Cells.Item(56, 10).ClearContents()
$ws.Cells.Item(56, 11).ClearContents()
$ws.Cells.Item(56, 12).ClearContents()
$ws.Cells.Item(56, 13).ClearContents()
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(57, 5).Value = 0
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 0
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 13).Value = 0
$ws.Cells.Item(58, 4).Value = 1000
$ws.Cells.Item(58, 5).Value = 1000
$ws.Cells.Item(58, 6).Value = 1000
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(58, 8).Value = 337000
$ws.Cells.Item(58, 9).Value = 336600
$ws.Cells.Item(58, 10).Value = 735900
$ws.Cells.Item(58, 11).Value = 735500
$ws.Cells.Item(58, 12).Value = "NA"
$ws.Cells.Item(58, 13).Value = "NA"
$ws.Cells.Item(59, 4).Value = 73049000
$ws.Cells.Item(59, 5).Value = 83034000
$ws.Cells.Item(59, 6).Value = 80599000
$ws.Cells.Item(59, 7).Value = 79668000
$ws.Cells.Item(59, 8).Value = 79471000
$ws.Cells.Item(59, 9).Value = 111603200
$ws.Cells.Item(59, 10).Value = 105725300
$ws.Cells.Item(59, 11).Value = 103261300
$ws.Cells.Item(59, 12).Value = 67883100
$ws.Cells.Item(59, 13).Value = 101627000
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(60, 5).Value = 0
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 8).Value = 0
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 11).Value = 0
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 13).Value = 0
$ws.Cells.Item(61, 4).Value = 3136000
$ws.Cells.Item(61, 5).Value = 3459000
$ws.Cells.Item(61, 6).Value = 3458000
$ws.Cells.Item(61, 7).Value = 3458000
$ws.Cells.Item(61, 8).Value = 3123000
$ws.Cells.Item(61, 9).Value = 3122200
$ws.Cells.Item(61, 10).Value = 2726200
$ws.Cells.Item(61, 11).Value = 2725700
$ws.Cells.Item(61, 12).Value = 3550000
$ws.Cells.Item(61, 13).Value = 3548500
$ws.Cells.Item(62, 4).Value = 551000
$ws.Cells.Item(62, 5).Value = 516000
$ws.Cells.Item(62, 6).Value = 527000
$ws.Cells.Item(62, 7).Value = 56998000
$ws.Cells.Item(62, 8).Value = 58827000
$ws.Cells.Item(62, 9).Value = 542200
$ws.Cells.Item(62, 10).Value = 636800
$ws.Cells.Item(62, 11).Value = 656000
$ws.Cells.Item(62, 12).Value = 60250000
$ws.Cells.Item(62, 13).Value = 631400
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = 0
$ws.Cells.Item(66, 4).Value = 146469000
$ws.Cells.Item(66, 5).Value = 156513000
$ws.Cells.Item(66, 6).Value = 154280000
$ws.Cells.Item(66, 7).Value = 210446000
$ws.Cells.Item(66, 8).Value = 212523000
$ws.Cells.Item(66, 9).Value = 212990900
$ws.Cells.Item(66, 10).Value = 206486600
$ws.Cells.Item(66, 11).Value = 204135100
$ws.Cells.Item(66, 12).Value = 201590000
$ws.Cells.Item(66, 13).Value = 203416200
$ws.Cells.Item(67, 4).ClearContents()
$ws.Cells.Item(67, 5).ClearContents()
$ws.Cells.Item(67, 6).ClearContents()
$ws.Cells.Item(67, 7).ClearContents()
$ws.Cells.Item(67, 8).ClearContents()
$ws.Cells.Item(67, 9).ClearContents()
$ws.Cells.Item(67, 10).ClearContents()
$ws.Cells.Item(67, 11).ClearContents()
$ws.Cells.Item(67, 12).ClearContents()
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 0
$ws.Cells.Item(72, 4).Value = -11732000
$ws.Cells.Item(72, 5).Value = -11853000
$ws.Cells.Item(72, 6).Value = -11995000
$ws.Cells.Item(72, 7).Value = -12161000
$ws.Cells.Item(72, 8).Value = -12719000
$ws.Cells.Item(72, 9).Value = -9655600
$ws.Cells.Item(72, 10).Value = -9804200
$ws.Cells.Item(72, 11).Value = -9971400
$ws.Cells.Item(72, 12).Value = -9742000
$ws.Cells.Item(72, 13).Value = -9310300
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).Value = 0
$ws.Cells.Item(76, 4).Value = 8213000
$ws.Cells.Item(76, 5).Value = 8523000
$ws.Cells.Item(76, 6).Value = 8460000
$ws.Cells.Item(76, 7).Value = 9378000
$ws.Cells.Item(76, 8).Value = 10009000
$ws.Cells.Item(76, 9).Value = 13653000
$ws.Cells.Item(76, 10).Value = 13352200
$ws.Cells.Item(76, 11).Value = 12890400
$ws.Cells.Item(76, 12).Value = 12995000
$ws.Cells.Item(76, 13).Value = 15206000
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(80, 6).Value = 43281
$ws.Cells.Item(80, 7).Value = 43190
$ws.Cells.Item(80, 8).Value = 43100
$ws.Cells.Item(80, 9).Value = 43008
$ws.Cells.Item(80, 10).Value = 42916
$ws.Cells.Item(80, 11).Value = 42825
$ws.Cells.Item(80, 12).Value = 42735
$ws.Cells.Item(80, 13).Value = 42643
$ws.Cells.Item(81, 4).Value = 121000
$ws.Cells.Item(81, 5).Value = 142000
$ws.Cells.Item(81, 6).Value = 166000
$ws.Cells.Item(81, 7).Value = 446000
$ws.Cells.Item(81, 8).Value = -3164300
$ws.Cells.Item(81, 9).Value = 149000
$ws.Cells.Item(81, 10).Value = 167000
$ws.Cells.Item(81, 11).Value = -143000
$ws.Cells.Item(81, 12).Value = -327000
$ws.Cells.Item(81, 13).Value = -248100
$ws.Cells.Item(82, 4).ClearContents()
$ws.Cells.Item(82, 5).ClearContents()
$ws.Cells.Item(82, 6).ClearContents()
$ws.Cells.Item(82, 7).ClearContents()
$ws.Cells.Item(82, 8).ClearContents()
$ws.Cells.Item(82, 9).ClearContents()
$ws.Cells.Item(82, 10).ClearContents()
$ws.Cells.Item(82, 11).ClearContents()
$ws.Cells.Item(82, 12).ClearContents()
$ws.Cells.Item(82, 13).ClearContents()
$ws.Cells.Item(83, 4).Value = 0
$ws.Cells.Item(83, 5).Value = 0
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = 0
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = 0
$ws.Cells.Item(89, 4).Value = 234000
$ws.Cells.Item(89, 5).Value = 580000
$ws.Cells.Item(89, 6).Value = 653000
$ws.Cells.Item(89, 7).Value = 401000
$ws.Cells.Item(89, 8).Value = 456000
$ws.Cells.Item(89, 9).Value = 570000
$ws.Cells.Item(89, 10).Value = 605800
$ws.Cells.Item(89, 11).Value = -49800
$ws.Cells.Item(89, 12).Value = 3591000
$ws.Cells.Item(89, 13).Value = 819700
$ws.Cells.Item(90, 4).ClearContents()
$ws.Cells.Item(90, 5).ClearContents()
$ws.Cells.Item(90, 6).ClearContents()
$ws.Cells.Item(90, 7).ClearContents()
$ws.Cells.Item(90, 8).ClearContents()
$ws.Cells.Item(90, 9).ClearContents()
$ws.Cells.Item(90, 10).ClearContents()
$ws.Cells.Item(90, 11).ClearContents()
$ws.Cells.Item(90, 12).ClearContents()
$ws.Cells.Item(90, 13).ClearContents()
$ws.Cells.Item(91, 4).Value = "NA"
$ws.Cells.Item(91, 5).Value = "NA"
$ws.Cells.Item(91, 6).Value = "NA"
$ws.Cells.Item(91, 7).Value = "NA"
$ws.Cells.Item(91, 8).Value = "NA"
$ws.Cells.Item(91, 9).Value = -9700
$ws.Cells.Item(91, 10).Value = -15700
$ws.Cells.Item(91, 11).Value = -10400
$ws.Cells.Item(91, 12).Value = -66700
$ws.Cells.Item(91, 13).Value = -25800
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 0
$ws.Cells.Item(94, 4).Value = -333000
$ws.Cells.Item(94, 5).Value = -36000
$ws.Cells.Item(94, 6).Value = -149000
$ws.Cells.Item(94, 7).Value = 236000
$ws.Cells.Item(94, 8).Value = -145000
$ws.Cells.Item(94, 9).Value = -1581300
$ws.Cells.Item(94, 10).Value = -892500
$ws.Cells.Item(94, 11).Value = 190800
$ws.Cells.Item(94, 12).Value = -3683000
$ws.Cells.Item(94, 13).Value = -1924000
$ws.Cells.Item(95, 4).ClearContents()
$ws.Cells.Item(95, 5).ClearContents()
$ws.Cells.Item(95, 6).ClearContents()
$ws.Cells.Item(95, 7).ClearContents()
$ws.Cells.Item(95, 8).ClearContents()
$ws.Cells.Item(95, 9).ClearContents()
$ws.Cells.Item(95, 10).ClearContents()
$ws.Cells.Item(95, 11).ClearContents()
$ws.Cells.Item(95, 12).ClearContents()
$ws.Cells.Item(95, 13).ClearContents()
$ws.Cells.Item(96, 4).Value = -1000
$ws.Cells.Item(96, 5).Value = -2000
$ws.Cells.Item(96, 6).Value = -1000
$ws.Cells.Item(96, 7).Value = -2000
$ws.Cells.Item(96, 8).Value = -2500
$ws.Cells.Item(96, 9).Value = -1700
$ws.Cells.Item(96, 10).Value = -1900
$ws.Cells.Item(96, 11).Value = -1900
$ws.Cells.Item(96, 12).Value = -8000
$ws.Cells.Item(96, 13).Value = -2000
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 0
$ws.Cells.Item(100, 4).Value = -152000
$ws.Cells.Item(100, 5).Value = -289000
$ws.Cells.Item(100, 6).Value = -926000
$ws.Cells.Item(100, 7).Value = -397000
$ws.Cells.Item(100, 8).Value = -562400
$ws.Cells.Item(100, 9).Value = 414200
$ws.Cells.Item(100, 10).Value = 556600
$ws.Cells.Item(100, 11).Value = -753400
$ws.Cells.Item(100, 12).Value = 495400
$ws.Cells.Item(100, 13).Value = 71300
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(101, 13).Value = 0
$ws.Cells.Item(102, 4).Value = -251000
$ws.Cells.Item(102, 5).Value = 255000
$ws.Cells.Item(102, 6).Value = -422000
$ws.Cells.Item(102, 7).Value = 240000
$ws.Cells.Item(102, 8).Value = -251200
$ws.Cells.Item(102, 9).Value = -600500
$ws.Cells.Item(102, 10).Value = 269100
$ws.Cells.Item(102, 11).Value = -612400
$ws.Cells.Item(102, 12).Value = 398000
$ws.Cells.Item(102, 13).Value = -1033000
